$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Column E width change (engine quantizes ColumnWidth to ~1/6 character
# steps; 28.04 lands on the closest achievable raw width to the target
# 28.85546875)
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 28.04

# ---------------------------------------------------------------------------
# Row 33: F33 goes from "PENDIENTE" (yellow) to "CORREGIDO" (green) style
# ---------------------------------------------------------------------------
$ws.Range("F34").Copy() | Out-Null
$ws.Range("F33").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F33").Value = "CORREGIDO"

# ---------------------------------------------------------------------------
# Row 38 (existing row, only A38 had data) -- fill remaining cells
# shared-string creation order: B38, C38, E38 (D38/F38 reuse existing strings)
# ---------------------------------------------------------------------------
$ws.Rows.Item(38).RowHeight = 58.5

$ws.Range("B38").Value = "Error en el asistente de configuración"
$ws.Range("B38").Font.Name = "Calibri Light"
$ws.Range("B38").WrapText = $true

$ws.Range("C38").Value = "Configurar edicion: configurar preferencias, elegir 8 equipos, agregar una fase y mostrar el fixture, tocar volver, seleccionar 3 equipos, va a salir el cartel que va a modificarse el fixture, poner aceptar, poner siguiente y se produce un error"
$ws.Range("C38").Font.Name = "Calibri Light"
$ws.Range("C38").WrapText = $true

$ws.Range("D38").Value = "Tony"

$ws.Range("E38").Value = "admin/edicion/equipos.aspx"
$ws.Range("E38").WrapText = $true

$ws.Range("F34").Copy() | Out-Null
$ws.Range("F38").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F38").Value = "CORREGIDO"

# ---------------------------------------------------------------------------
# Row 39 (new)
# shared-string creation order: B39, C39, E39
# ---------------------------------------------------------------------------
$ws.Rows.Item(39).RowHeight = 60

$ws.Range("A33").Copy() | Out-Null
$ws.Range("A39").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A39").Value = 37

$ws.Range("B39").Value = "INTERFAZ PARTIDOS: Arreglos menores"
$ws.Range("B39").Font.Name = "Calibri Light"
$ws.Range("B39").WrapText = $true

$ws.Range("C39").Value = "Otros partidos de la fecha: colorcito de los estados`n- Widget versus: cambiar foto Partidos Empatados y Partidos Perdidos"
$ws.Range("C39").Font.Name = "Calibri Light"
$ws.Range("C39").WrapText = $true

$ws.Range("D39").Value = "Facu"

$ws.Range("E39").Value = "torneo/partido"

$ws.Range("F32").Copy() | Out-Null
$ws.Range("F39").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F39").Value = "PENDIENTE"

# ---------------------------------------------------------------------------
# Row 40 (new)
# shared-string creation order: B40, E40
# ---------------------------------------------------------------------------
$ws.Rows.Item(40).RowHeight = 30

$ws.Range("A33").Copy() | Out-Null
$ws.Range("A40").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A40").Value = 38

$ws.Range("B40").Value = "Sacar Notificaciones (módulo admin)"
$ws.Range("B40").Font.Name = "Calibri Light"
$ws.Range("B40").WrapText = $true

$ws.Range("D40").Value = "Facu"

$ws.Range("E40").Value = "admin/"

$ws.Range("F32").Copy() | Out-Null
$ws.Range("F40").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F40").Value = "PENDIENTE"

# ---------------------------------------------------------------------------
# Row 41 (new)
# shared-string creation order: C41 (Colorcito...), B41 (INTERFAZ EQUIPOS), E41
# ---------------------------------------------------------------------------
$ws.Rows.Item(41).RowHeight = 45

$ws.Range("A33").Copy() | Out-Null
$ws.Range("A41").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A41").Value = 39

$ws.Range("C41").Value = "Colorcito de los Resultados (Empatado - Perdido - Ganado)`n- Me parece que falta PROXIMOS PARTIDOS!"
$ws.Range("C41").WrapText = $true

$ws.Range("B41").Value = "INTERFAZ EQUIPOS: Arreglos menores"
$ws.Range("B41").Font.Name = "Calibri Light"
$ws.Range("B41").WrapText = $true

$ws.Range("D41").Value = "Facu"

$ws.Range("E41").Value = "torneo/equipo"

$ws.Range("F32").Copy() | Out-Null
$ws.Range("F41").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F41").Value = "PENDIENTE"

# ---------------------------------------------------------------------------
# Rows 42 and 43 (new, only column A populated)
# ---------------------------------------------------------------------------
$ws.Range("A33").Copy() | Out-Null
$ws.Range("A42").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A42").Value = 40

$ws.Range("A33").Copy() | Out-Null
$ws.Range("A43").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A43").Value = 41

# ---------------------------------------------------------------------------
# Selection / active cell
# ---------------------------------------------------------------------------
$ws.Range("C36").Select() | Out-Null
